$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 18..48 (col A = station name, col B = terminal name, col C = datetime serial)
$rows = @(
    ,@("长沙特来电飞狐四方坪西区充电站", "101号直流", 45935.0465625)
    ,@("长沙特来电飞狐四方坪南区充电站", "406号直流", 45943.020914351851)
    ,@("长沙特来电飞狐四方坪南区充电站", "201号直流", 45944.074282407404)
    ,@("长沙特来电飞狐四方坪南区充电站", "203号直流", 45944.228055555555)
    ,@("长沙特来电飞狐四方坪西区充电站", "502号直流", 45944.54005787037)
    ,@("长沙特来电飞狐四方坪南区充电站", "306号直流", 45944.674907407411)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "110号直流", 45944.829027777778)
    ,@("长沙特来电飞狐四方坪南区充电站", "902号直流", 45945.197604166664)
    ,@("长沙特来电飞狐四方坪西区充电站", "702号直流", 45945.22991898148)
    ,@("长沙特来电飞狐四方坪西区充电站", "505号直流", 45945.273923611108)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "305号直流", 45945.510601851849)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "104号直流", 45945.517974537041)
    ,@("长沙特来电飞狐四方坪东区充电站", "904号直流", 45945.536215277774)
    ,@("长沙特来电飞狐四方坪南区充电站", "104号直流", 45945.548738425925)
    ,@("长沙特来电飞狐四方坪西区充电站", "B03号直流", 45945.562523148146)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "301号直流", 45945.567685185182)
    ,@("长沙特来电飞狐四方坪西区充电站", "801号直流", 45945.568148148152)
    ,@("长沙特来电飞狐四方坪西区充电站", "B01号直流", 45945.572557870371)
    ,@("长沙特来电飞狐四方坪西区充电站", "804号直流", 45945.583333333336)
    ,@("长沙特来电飞狐四方坪东区充电站", "102号直流", 45945.589571759258)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "112号直流", 45945.592013888891)
    ,@("长沙特来电飞狐四方坪西区充电站", "602号直流", 45945.598043981481)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "109号直流", 45945.603726851848)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "208号直流", 45945.617673611108)
    ,@("长沙特来电飞狐四方坪南区充电站", "305号直流", 45945.6328125)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "106号直流", 45945.663530092592)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "108号直流", 45945.68378472222)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "101号直流", 45945.720682870371)
    ,@("长沙市开福区高岭香江国际城充电站建设项目", "107号直流", 45945.729618055557)
    ,@("长沙特来电飞狐四方坪南区充电站", "301号直流", 45945.730069444442)
    ,@("长沙特来电飞狐四方坪西区充电站", "903号直流", 45945.73505787037)
)

$startRow = 18
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $entry = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
}

# Update the view state: scrolled position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("H24").Select()
